$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44680
$ws.Range("J2").Value = 20
$ws.Range("D3").Value = 44259
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = 4000
$ws.Range("P3").Value = 4000
$ws.Range("D4").Value = 44749
$ws.Range("J4").Value = 65
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("D5").Value = 45162
$ws.Range("D6").Value = 44957
$ws.Range("J6").Value = 20
$ws.Range("D7").Value = 44656
$ws.Range("J7").Value = 85
$ws.Range("D8").Value = 45159
$ws.Range("J8").Value = 75
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = 5000
$ws.Range("P8").Value = 5000
$ws.Range("D10").Value = 44301
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("P10").Value = 3000
$ws.Range("D11").Value = 44312
$ws.Range("J11").Value = 50
$ws.Range("D12").Value = 44390
$ws.Range("J12").Value = 55
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 6000
$ws.Range("P12").Value = 6000
$ws.Range("D13").Value = 44959
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 5000
$ws.Range("P13").Value = 5000
$ws.Range("D14").Value = 44365
$ws.Range("J14").Value = 55
$ws.Range("D15").Value = 44316
$ws.Range("J15").Value = 20
$ws.Range("D16").Value = 44176
$ws.Range("J16").Value = 10
$ws.Range("D17").Value = 44781
$ws.Range("J17").Value = 40
$ws.Range("D18").Value = 44649
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 5000
$ws.Range("P18").Value = 5000
$ws.Range("D19").Value = 44313
$ws.Range("J19").Value = 20
$ws.Range("L19").Value = 4000
$ws.Range("M19").Value = 4000
$ws.Range("P19").Value = 4000
$ws.Range("D20").Value = 45169
$ws.Range("K20").Value = 4000
$ws.Range("M20").Value = 4600
$ws.Range("P20").Value = 4600
$ws.Range("D21").Value = 44509
$ws.Range("J21").Value = 20
$ws.Range("D22").Value = 44777
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = 5000
$ws.Range("P22").Value = 5000
$ws.Range("D23").Value = 44291
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = 4000
$ws.Range("L23").Value = 4000
$ws.Range("M23").Value = 4000
$ws.Range("P23").Value = 4000
$ws.Range("D24").Value = 44280
$ws.Range("J24").Value = 55
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = 4000
$ws.Range("P24").Value = 4000
$ws.Range("D25").Value = 44956
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 5000
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 5000
$ws.Range("P25").Value = 5000
$ws.Range("D26").Value = 45163
$ws.Range("J26").Value = 30
$ws.Range("D27").Value = 44504
$ws.Range("J27").Value = 55
$ws.Range("D28").Value = 44679
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = 5000
$ws.Range("P28").Value = 5000
$ws.Range("D29").Value = 44497
$ws.Range("K29").Value = 4000
$ws.Range("L29").Value = 4000
$ws.Range("M29").Value = 4000
$ws.Range("P29").Value = 4000
$ws.Range("D30").Value = 44966
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = 5000
$ws.Range("P30").Value = 5000
$ws.Range("D31").Value = 44498
$ws.Range("D32").Value = 44315
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = 4000
$ws.Range("P32").Value = 4000
$ws.Range("D33").Value = 45194
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 6000
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = 6000
$ws.Range("P33").Value = 6000
